# Swap the order of names in the "Recorded By" (column G) cells that list
# both "dnasr281@gmail.com" and "System" as recorders, changing
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
# for every row on the active worksheet where this occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is the "Recorded By" column.
$col = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
